$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Normalize the reversed-order injury-pair labels in the Station columns (K:R) ---
# (the lab re-coded these so the two-ligament labels always list the "primary" ligament first)
$ws.Range("R4").Value  = "ACL, MCL"
$ws.Range("N5").Value  = "ACL, LCL"
$ws.Range("O5").Value  = "ACL, LCL"
$ws.Range("Q5").Value  = "ACL, MCL"
$ws.Range("R5").Value  = "ACL, LCL"
$ws.Range("K6").Value  = "MCL, PCL"
$ws.Range("M6").Value  = "LCL, PCL"
$ws.Range("R6").Value  = "MCL, PCL"
$ws.Range("K7").Value  = "MCL, PCL"
$ws.Range("O9").Value  = "MCL, PCL"

# --- 2) Fold the separate "Stand Alone" skills table (A13:F21) into the main table ---
# as five new columns (BC:BG), keyed by Subject, then remove the old standalone table.
$ws.Range("BC1").Value = "GainedNew"
$ws.Range("BD1").Value = "Improved"
$ws.Range("BE1").Value = "Effective"
$ws.Range("BF1").Value = "Retained"
$ws.Range("BG1").Value = "Transferable"

$ws.Range("BC2").Value = 10
$ws.Range("BD2").Value = 10
$ws.Range("BE2").Value = 10
$ws.Range("BF2").Value = 10
$ws.Range("BG2").Value = 8

$ws.Range("BC3").Value = 10
$ws.Range("BD3").Value = 9
$ws.Range("BE3").Value = "N/A"
$ws.Range("BF3").Value = "N/A"
$ws.Range("BG3").Value = "N/A"

$ws.Range("BC4").Value = 10
$ws.Range("BD4").Value = 10
$ws.Range("BE4").Value = 10
$ws.Range("BF4").Value = 10
$ws.Range("BG4").Value = 7

$ws.Range("BC5").Value = 10
$ws.Range("BD5").Value = 10
$ws.Range("BE5").Value = 8
$ws.Range("BF5").Value = 8
$ws.Range("BG5").Value = 8

$ws.Range("BC6").Value = 9
$ws.Range("BD6").Value = 7
$ws.Range("BE6").Value = 9
$ws.Range("BF6").Value = 9
$ws.Range("BG6").Value = 9

$ws.Range("BC7").Value = 10
$ws.Range("BD7").Value = 10
$ws.Range("BE7").Value = 10
$ws.Range("BF7").Value = 10
$ws.Range("BG7").Value = 10

$ws.Range("BC8").Value = 10
$ws.Range("BD8").Value = 10
$ws.Range("BE8").Value = 10
$ws.Range("BF8").Value = 10
$ws.Range("BG8").Value = 10

$ws.Range("BC9").Value = 9
$ws.Range("BD9").Value = 7
$ws.Range("BE9").Value = 8
$ws.Range("BF9").Value = 8
$ws.Range("BG9").Value = 7

# Remove the now-redundant standalone table entirely (values + formatting).
$ws.Range("A13:F21").Clear()

# A few cells below keep the carried-over "no value, header-style font" formatting
# that was left behind once the old table's rows were cleared out.
$ws.Range("H1").Copy()
$ws.Range("B20:B22").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 3) Leave the selection where the editor last left it ---
$ws.Range("R12").Select()
